$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of row 2 into new row 3 (add import value validation row)
$ws.Range("A2:L2").Copy()
$ws.Range("A3:L3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Set values for the new row 3
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = "'0200"
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = "Nacht"
$ws.Range("H3").Value = 35
$ws.Range("I3").Value = 2344
$ws.Range("J3").Value = 1.5678000000000001
$ws.Range("K3").Value = 42.548999999999999
$ws.Range("L3").Value = 1.5

# Update active selection to J3 as per the saved workbook view
$ws.Range("J3").Select()
